# Fix: "avg_degree" metric in column E was computed as (edge_count / vertex_count)
# for the last set of rows instead of the correct (2 * edge_count / vertex_count),
# which made the metric understated (effectively half of what it should be).
# Correct it by doubling the existing values in E22:E37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 22; $row -le 37; $row++) {
    $cell = $ws.Cells.Item($row, 5)  # column E = 5
    $cell.Value = $cell.Value2 * 2
}
